{"js": "// Split the single \"Bibliografia\" run into multiple runs separated by\n// <w:br/> line breaks, one before each numbered reference ([1]..[5]),\n// with an extra blank line (double break) right before [1].\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the bibliography paragraph robustly (rather than hard-coding index).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text;\n  if (t.indexOf(\"A bibliografia ser\u00e1 recomendada\") !== -1 && t.indexOf(\"[1]\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  throw new Error(\"Bibliography paragraph not found\");\n}\n\n// Insert breaks in reverse order ([5] -> [2]) so earlier search hits are not\n// shifted by later insertions; finish with [1], which needs two breaks\n// (i.e. a blank line) immediately after \"mentoria.\".\nconst singleBreakMarkers = [\n  \"[5] Diretrizes Curriculares Nacionais\",\n  \"[4] Kaul, S. Triangulated Mentorship\",\n  \"[3] Mueller, S. Electronic mentoring\",\n  \"[2] Zachary, L. J. The Mentor\",\n];\n\nfor (const marker of singleBreakMarkers) {\n  const results = target.search(marker, { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '\" + marker + \"', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"\\u000b\", Word.InsertLocation.before);\n  await context.sync();\n}\n\n{\n  const results = target.search(\"[1] Peddy, S. The art of mentoring\", { matchCase: true, matchWildcards: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error(\"Expected exactly one match for '[1] Peddy...', found \" + results.items.length);\n  }\n  results.items[0].insertText(\"\\u000b\\u000b\", Word.InsertLocation.before);\n  await context.sync();\n}\n", "ps1": "# Split the single \"Bibliografia\" run into multiple runs separated by\n# line breaks (<w:br/>), one before each numbered reference ([1]..[5]),\n# with an extra blank line (double break) right before [1].\n\n$d = $word.ActiveDocument\n\n$bib = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*bibliografia*\" -and $t -like \"*`[1`]*\" -and $t -like \"*Peddy*\") {\n        $bib = $p\n    }\n}\nif ($bib -eq $null) {\n    throw \"Bibliography paragraph not found\"\n}\n\n# Insert breaks from [5] down to [2] first so earlier matches in the text\n# aren't shifted by later insertions; [1] is handled last with two breaks\n# (a blank line) right after \"mentoria.\".\n$singleBreakMarkers = @(\n    \"[5] Diretrizes Curriculares Nacionais\",\n    \"[4] Kaul, S. Triangulated Mentorship\",\n    \"[3] Mueller, S. Electronic mentoring\",\n    \"[2] Zachary, L. J. The Mentor\"\n)\n\nforeach ($marker in $singleBreakMarkers) {\n    $rng = $bib.Range.Duplicate()\n    $found = $rng.Find.Execute($marker)\n    if (-not $found) {\n        throw \"Marker not found: $marker\"\n    }\n    $rng.Collapse(1)\n    $rng.InsertBefore([char]11)\n}\n\n$rng = $bib.Range.Duplicate()\n$found = $rng.Find.Execute(\"[1] Peddy, S. The art of mentoring\")\nif (-not $found) {\n    throw \"Marker not found: [1] Peddy...\"\n}\n$rng.Collapse(1)\n$rng.InsertBefore([string]([char]11) + [string]([char]11))\n"}
